# Automatic push - Update excel files
#
# The original edit refreshed the cached values that come from an external
# workbook link ("Historique d'achats.xlsx", referenced in formulas as
# "[2]<Sheet>!<Cell>", e.g. =[2]ETH!J4). Excel keeps a last-known-value
# cache for such links and updates it when the link is refreshed against
# the live source file. That source file is not available in this
# environment, so we reproduce the net effect of the refresh by writing
# the new cached numbers straight into the cells that pull the external
# data on sheet "Feuil1". Every other changed cell in the workbook (totals,
# ratios, the pie-chart cache, etc.) is a formula that depends on these
# cells, so it recalculates to the correct value automatically once Excel
# recomputes the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$newValues = @{
    "C12" = 1724.8158243988371   # [2]ETH!J4
    "C13" = 1301.6010433248623   # [2]BTC!J4
    "C14" = 476.07806721939266   # [2]SOL!J4
    "C17" = 241.7794083554943    # [2]BNB!J4
    "C20" = 60.845173936385557   # [2]AVAX!$J$4
    "C21" = 53.827562698611125   # [2]MATIC!$J$4
    "C22" = 55.731439807770421   # [2]NEAR!$J$4
    "C24" = 48.609356726245665   # [2]DOT!$J$4
    "C25" = 45.901047779934125   # [2]LUNC!J4
    "C27" = 43.805776956255514   # [2]ADA!$J$4
    "C28" = 48.79152517417927    # [2]MINA!$J$4
    "C29" = 23.738689379963137   # [2]TIA!$J$4
    "C30" = 25.236837662854541   # [2]APE!$J$4
    "C31" = 20.652338384945629   # [2]DYDX!$J$4
    "C32" = 17.993131105025743   # [2]UNI!$J$4
    "C33" = 17.967407331151506   # [2]LDO!$J$4
    "C34" = 16.385640811155884   # [2]XRP!$J$4
    "C35" = 14.442125046863765   # [2]SHIB!$J$4
    "C36" = 13.971891538579618   # [2]ICP!$J$4
    "C37" = 13.484146487787511   # [2]LINK!$J$4
    "C38" = 12.645374525688341   # [2]ATOM!$J$4
    "C39" = 11.336290612122465   # [2]LTC!$J$4
    "C40" = 11.604690661799797   # [2]ALGO!$J$4
    "C42" = 6.1403056472286819   # [2]EGLD!$J$4
    "C43" = 4.9886608432033732   # [2]DOGE!$J$4
    "C44" = 4.9018349848171612   # [2]LUNA!J4
    "C45" = 4.6531377859358001   # [2]GRT!$J$4
    "C46" = 3.4987492538727851   # [2]AMP!$J$4
    "C47" = 2.9212988888518234   # [2]ACE!$J$4
    "C48" = 2.8877486200486229   # [2]SEI!$J$4
    "C49" = 3.0893617676439189   # [2]SHPING!$J$4
    "C50" = 2.1099637931259538   # [2]KAVA!$J$4
    "C51" = 2.6767265945147103   # [2]POLIS!J4
    "C52" = 1.5850885760236557   # [2]MEME!$J$4
    "C53" = 1.188085806013488    # [2]TRX!$J$4
    "C54" = 0.50898428667120044  # [2]ATLAS!O47
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

$excel.CalculateFull()
